$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 250, shifting existing rows 250:301 down to 251:302.
$ws.Rows("250:250").Insert()

# Populate the newly inserted row 250 with the new observation.
$ws.Cells.Item(250, 1).Value = 5
$ws.Cells.Item(250, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(250, 3).Value = "Maule"
$ws.Cells.Item(250, 4).Value = 44637
$ws.Cells.Item(250, 5).Value = 7
$ws.Cells.Item(250, 6).Value = 100114013
$ws.Cells.Item(250, 7).Value = "Zanahoria"
$ws.Cells.Item(250, 8).Value = "Sin especificar"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 500
$ws.Cells.Item(250, 11).Value = 7000
$ws.Cells.Item(250, 12).Value = 7000
$ws.Cells.Item(250, 13).Value = 7000
$ws.Cells.Item(250, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(250, 15).Value = "Región de Ñuble"
$ws.Cells.Item(250, 16).Value = 350
$ws.Cells.Item(250, 17).Value = 20
$ws.Cells.Item(250, 18).Value = "Hortaliza"
